$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 195.75
$ws.Range("I9").Value = 198.71428
$ws.Range("J9").Value = 175
$ws.Range("K9").Value = 198.71428
$ws.Range("L9").Value = 175
$ws.Range("M9").Value = -29.71428
$ws.Range("N9").Value = -513
$ws.Range("H15").Value = 1687.2295
$ws.Range("I15").Value = 1687.2295
$ws.Range("K15").Value = 5061.6885
$ws.Range("M15").Value = -4892.6885
$ws.Range("H121").Value = 5014.375
$ws.Range("I121").Value = 1900
$ws.Range("J121").Value = 6883
$ws.Range("K121").Value = 5700
$ws.Range("L121").Value = 20649
$ws.Range("M121").Value = -3953
$ws.Range("N121").Value = -24143
$ws.Range("H137").Value = 2568.9565
$ws.Range("I137").Value = 2539.0667
$ws.Range("K137").Value = 7617.2001
$ws.Range("M137").Value = -5067.2001
$ws.Range("H141").Value = 4948.864
$ws.Range("I141").Value = 4037.5
$ws.Range("J141").Value = 6042.5
$ws.Range("K141").Value = 12112.5
$ws.Range("L141").Value = 18127.5
$ws.Range("M141").Value = -6932.5
$ws.Range("N141").Value = -28487.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 22500
$ws.Range("I31").Value = 22500
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 22500
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -22206
$ws.Range("N31").ClearContents()
$ws.Range("H45").Value = 1685278.4
$ws.Range("I45").Value = 2332647.8
$ws.Range("J45").Value = 2117.6
$ws.Range("K45").Value = 2332647.8
$ws.Range("L45").Value = 2117.6
$ws.Range("M45").Value = -2332270.8
$ws.Range("N45").Value = -2871.6
$ws.Range("H74").Value = 5290.7417
$ws.Range("I74").Value = 2174.2593
$ws.Range("J74").Value = 26327
$ws.Range("K74").Value = 2174.2593
$ws.Range("L74").Value = 26327
$ws.Range("M74").Value = -1300.2593
$ws.Range("N74").Value = -28075
$ws.Range("H77").Value = 5290.7417
$ws.Range("I77").Value = 2174.2593
$ws.Range("J77").Value = 26327
$ws.Range("K77").Value = 10871.2965
$ws.Range("L77").Value = 131635
$ws.Range("M77").Value = -6503.2965
$ws.Range("N77").Value = -140371
$ws.Range("H130").Value = 70400
$ws.Range("J130").Value = 70400
$ws.Range("L130").Value = 70400
$ws.Range("N130").Value = -80440
$ws.Range("H132").Value = 1661.8572
$ws.Range("I132").Value = 1358.3226
$ws.Range("J132").Value = 2517.2727
$ws.Range("K132").Value = 4074.9678
$ws.Range("L132").Value = 7551.8181
$ws.Range("M132").Value = -1544.9678
$ws.Range("N132").Value = -12611.8181
$ws.Range("H134").Value = 37476.332
$ws.Range("J134").Value = 37476.332
$ws.Range("L134").Value = 37476.332
$ws.Range("N134").Value = -47616.332

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H86").Value = 5748915
$ws.Range("I86").Value = 6537712.5
$ws.Range("J86").Value = 1959.1428
$ws.Range("K86").Value = 6537712.5
$ws.Range("L86").Value = 1959.1428
$ws.Range("M86").Value = -6536589.5
$ws.Range("N86").Value = -4205.1428
$ws.Range("H89").Value = 5748915
$ws.Range("I89").Value = 6537712.5
$ws.Range("J89").Value = 1959.1428
$ws.Range("K89").Value = 32688562.5
$ws.Range("L89").Value = 9795.714
$ws.Range("M89").Value = -32682946.5
$ws.Range("N89").Value = -21027.714
$ws.Range("H134").Value = 38862.703
$ws.Range("I134").Value = 1771.5294
$ws.Range("K134").Value = 5314.5882
$ws.Range("M134").Value = -2779.5882

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2842708.5
$ws.Range("I58").Value = 4133403.2
$ws.Range("K58").Value = 4133403.2
$ws.Range("M58").Value = -4133200.2
$ws.Range("H63").Value = 40271
$ws.Range("J63").Value = 40271
$ws.Range("L63").Value = 40271
$ws.Range("N63").Value = -41643
$ws.Range("H66").Value = 40271
$ws.Range("J66").Value = 40271
$ws.Range("L66").Value = 120813
$ws.Range("N66").Value = -127677
$ws.Range("H132").Value = 2250.7036
$ws.Range("I132").Value = 2122.0303
$ws.Range("J132").Value = 2452.9048
$ws.Range("K132").Value = 6366.090899999999
$ws.Range("L132").Value = 7358.714399999999
$ws.Range("M132").Value = -3836.090899999999
$ws.Range("N132").Value = -12418.7144
$ws.Range("H134").Value = 2442.3845
$ws.Range("I134").Value = 2084.2334
$ws.Range("J134").Value = 3636.2222
$ws.Range("K134").Value = 6252.7002
$ws.Range("L134").Value = 10908.6666
$ws.Range("M134").Value = -3717.7002
$ws.Range("N134").Value = -15978.6666
$ws.Range("H136").Value = 2842708.5
$ws.Range("I136").Value = 4133403.2
$ws.Range("K136").Value = 12400209.6
$ws.Range("M136").Value = -12397659.6

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1653.1
$ws.Range("H134").Value = 3911.6
$ws.Range("I134").Value = 3328.9473
$ws.Range("J134").Value = 4918
$ws.Range("K134").Value = 9986.841899999999
$ws.Range("L134").Value = 14754
$ws.Range("M134").Value = -4916.841899999999
$ws.Range("N134").Value = -24894
$ws.Range("H137").Value = 20971.45
$ws.Range("J137").Value = 29900.6
$ws.Range("L137").Value = 89701.79999999999
$ws.Range("N137").Value = -99901.79999999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7189.364
$ws.Range("I132").Value = 2195.4
$ws.Range("K132").Value = 6586.200000000001
$ws.Range("M132").Value = -4056.200000000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6069.25
$ws.Range("I7").Value = 4618.909
$ws.Range("J7").Value = 9260
$ws.Range("K7").Value = 4618.909
$ws.Range("L7").Value = 9260
$ws.Range("M7").Value = -4506.909
$ws.Range("N7").Value = -9484
$ws.Range("H22").Value = 734
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 901
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 901
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1491
$ws.Range("H27").Value = 734
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 901
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 901
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -1115
$ws.Range("H40").Value = 3753.7693
$ws.Range("I40").Value = 3436.2727
$ws.Range("K40").Value = 3436.2727
$ws.Range("M40").Value = -3300.2727
$ws.Range("H62").Value = 40249
$ws.Range("J62").Value = 40249
$ws.Range("L62").Value = 40249
$ws.Range("N62").Value = -41497
$ws.Range("H65").Value = 40249
$ws.Range("J65").Value = 40249
$ws.Range("L65").Value = 120747
$ws.Range("N65").Value = -126987
$ws.Range("H122").Value = 7756.189
$ws.Range("I122").Value = 7303.2173
$ws.Range("J122").Value = 8500.357
$ws.Range("K122").Value = 21909.6519
$ws.Range("L122").Value = 25501.071
$ws.Range("M122").Value = -19459.6519
$ws.Range("N122").Value = -30401.071
$ws.Range("H126").Value = 6069.25
$ws.Range("I126").Value = 4618.909
$ws.Range("J126").Value = 9260
$ws.Range("K126").Value = 13856.727
$ws.Range("L126").Value = 27780
$ws.Range("M126").Value = -11386.727
$ws.Range("N126").Value = -32720

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 38967.668
$ws.Range("J80").Value = 38967.668
$ws.Range("L80").Value = 38967.668
$ws.Range("N80").Value = -40963.668
$ws.Range("H83").Value = 38967.668
$ws.Range("J83").Value = 38967.668
$ws.Range("L83").Value = 116903.004
$ws.Range("N83").Value = -126887.004
$ws.Range("H113").Value = 683.63336
$ws.Range("I113").Value = 363.6842
$ws.Range("J113").Value = 1236.2727
$ws.Range("K113").Value = 1091.0526
$ws.Range("L113").Value = 3708.8181
$ws.Range("M113").Value = 1078.9474
$ws.Range("N113").Value = -8048.8181
$ws.Range("H136").Value = 5736.0264
$ws.Range("I136").Value = 1948.3636
$ws.Range("J136").Value = 10944.0625
$ws.Range("K136").Value = 5845.0908
$ws.Range("L136").Value = 32832.1875
$ws.Range("M136").Value = -3295.0908
$ws.Range("N136").Value = -37932.1875

Write-Host "Applied all Pandaemonium Profits updates"